$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change A2 from "Joao" to "Lucas" (B2 stays 0)
$ws.Range("A2").Value = "Lucas"
$ws.Range("B2").Value = 0

# Add new rows: Breno (row 3) and I (row 4), each with 0 devices
$ws.Range("A3").Value = "Breno"
$ws.Range("B3").Value = 0

$ws.Range("A4").Value = "I"
$ws.Range("B4").Value = 0
